# edit.ps1 -- apply the skeleton.docx revision described by the commit diff:
#   1) "... there is just, probable and reasonable cause to believe that there
#      is now: " -> "... there is just, probable and reasonable cause to
#      believe that: "   (drop "there is now")
#   2) the {{ TRAININGEXPERIENCE }} placeholder token is renamed to
#      {{ T_AND_E }}
#   3) a stray empty paragraph right before the
#      "Together with other fruits, instrumentalities, and evidence of the
#      crime(s) of:" heading paragraph is removed.

$d = $word.ActiveDocument

# --- 1) "... there is now: " -> "... : " -------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "there is just, probable and reasonable cause to believe that there is now:"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "there is just, probable and reasonable cause to believe that:"
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null

# --- 2) rename the {{ TRAININGEXPERIENCE }} placeholder to {{ T_AND_E }} ---
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "TRAININGEXPERIENCE"
$find2.Replacement.ClearFormatting()
$find2.Replacement.Text = "T_AND_E"
$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null

# --- 3) delete the empty paragraph just before the "Together with other ---
# --- fruits, instrumentalities, and evidence of the crime(s) of:" line ----
$markerRange = $d.Content
$markerRange.Find.ClearFormatting()
$markerRange.Find.Execute("Together with other fruits, instrumentalities, and evidence of the crime(s) of:") | Out-Null
$headingParagraph = $markerRange.Paragraphs(1)
$emptyParagraph = $headingParagraph.Previous()
if ($emptyParagraph.Range.Text.Trim() -eq "") {
    $emptyParagraph.Range.Delete()
}
